$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44357
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 38000
$ws.Range("O2").Value = 38000
$ws.Range("P2").Value = 38000
$ws.Range("R2").Value = 'Perú'
$ws.Range("S2").Value = 2111
$ws.Range("D3").Value = 44629
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 35000
$ws.Range("O3").Value = 35000
$ws.Range("P3").Value = 35000
$ws.Range("S3").Value = 1944
$ws.Range("D4").Value = 44431
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 35000
$ws.Range("O4").Value = 35000
$ws.Range("P4").Value = 35000
$ws.Range("R4").Value = 'Región de Arica y Parinacota'
$ws.Range("S4").Value = 1944
$ws.Range("D5").Value = 44364
$ws.Range("M5").Value = 90
$ws.Range("N5").Value = 1700
$ws.Range("O5").Value = 1700
$ws.Range("P5").Value = 1700
$ws.Range("Q5").Value = '$/kilo'
$ws.Range("S5").Value = 1700
$ws.Range("T5").Value = 1
$ws.Range("D6").Value = 44294
$ws.Range("M6").Value = 15
$ws.Range("D7").Value = 44418
$ws.Range("M7").Value = 30
$ws.Range("D8").Value = 44264
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 40000
$ws.Range("O8").Value = 40000
$ws.Range("P8").Value = 40000
$ws.Range("S8").Value = 2222
$ws.Range("D9").Value = 44279
$ws.Range("M9").Value = 30
$ws.Range("O9").Value = 36000
$ws.Range("P9").Value = 35667
$ws.Range("S9").Value = 1982
$ws.Range("D10").Value = 44448
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 38000
$ws.Range("O10").Value = 38000
$ws.Range("P10").Value = 38000
$ws.Range("S10").Value = 2111
$ws.Range("D11").Value = 44379
$ws.Range("M11").Value = 10
$ws.Range("N11").Value = 30000
$ws.Range("O11").Value = 30000
$ws.Range("P11").Value = 30000
$ws.Range("S11").Value = 1667
$ws.Range("D12").Value = 44449
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 38000
$ws.Range("O12").Value = 38000
$ws.Range("P12").Value = 38000
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("S12").Value = 2111
$ws.Range("T12").Value = 18
$ws.Range("D13").Value = 44392
$ws.Range("M13").Value = 20
$ws.Range("R13").Value = 'Región de Arica y Parinacota'
$ws.Range("D14").Value = 44433
$ws.Range("M14").Value = 15
$ws.Range("D15").Value = 44377
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 40000
$ws.Range("O15").Value = 40000
$ws.Range("P15").Value = 40000
$ws.Range("R15").Value = 'Región de Arica y Parinacota'
$ws.Range("S15").Value = 2222
$ws.Range("D16").Value = 44369
$ws.Range("M16").Value = 5
$ws.Range("R16").Value = 'Perú'
$ws.Range("D17").Value = 44363
$ws.Range("M17").Value = 144
$ws.Range("N17").Value = 1700
$ws.Range("O17").Value = 1700
$ws.Range("P17").Value = 1700
$ws.Range("Q17").Value = '$/kilo'
$ws.Range("R17").Value = 'Región de Arica y Parinacota'
$ws.Range("S17").Value = 1700
$ws.Range("T17").Value = 1
$ws.Range("D18").Value = 44424
$ws.Range("M18").Value = 15
$ws.Range("D19").Value = 44442
$ws.Range("M19").Value = 15
$ws.Range("N19").Value = 35000
$ws.Range("O19").Value = 35000
$ws.Range("P19").Value = 35000
$ws.Range("R19").Value = 'Perú'
$ws.Range("S19").Value = 1944
$ws.Range("D20").Value = 44645
$ws.Range("M20").Value = 5
$ws.Range("D21").Value = 44405
$ws.Range("M21").Value = 10
$ws.Range("D22").Value = 44434
$ws.Range("M22").Value = 40
$ws.Range("D23").Value = 44634
$ws.Range("M23").Value = 30
$ws.Range("N23").Value = 45000
$ws.Range("O23").Value = 45000
$ws.Range("P23").Value = 45000
$ws.Range("S23").Value = 2500
$ws.Range("D24").Value = 44432
$ws.Range("M24").Value = 10
$ws.Range("N24").Value = 35000
$ws.Range("O24").Value = 35000
$ws.Range("P24").Value = 35000
$ws.Range("Q24").Value = '$/caja 18 kilos'
$ws.Range("R24").Value = 'Perú'
$ws.Range("S24").Value = 1944
$ws.Range("T24").Value = 18
$ws.Range("D25").Value = 44438
$ws.Range("M25").Value = 25
$ws.Range("N25").Value = 35000
$ws.Range("O25").Value = 35000
$ws.Range("P25").Value = 35000
$ws.Range("S25").Value = 1944
$ws.Range("D26").Value = 44435
$ws.Range("M26").Value = 10
$ws.Range("N26").Value = 35000
$ws.Range("O26").Value = 35000
$ws.Range("P26").Value = 35000
$ws.Range("R26").Value = 'Perú'
$ws.Range("S26").Value = 1944
$ws.Range("D27").Value = 44435
$ws.Range("M27").Value = 105
$ws.Range("R27").Value = 'Región de Arica y Parinacota'
